# Instantiate & Destroy Example.pptx
# - The 4th slide (p:sldId 288) contains a picture (cNvPr id="1044", name
#   "그림 145") that was re-saved by the original image-edit tool with:
#     * a "descr" (AlternativeText) pointing at the tool's temp-folder copy
#       of the source PNG
#     * a very slightly nudged position/size (it was nudged/"snapped" by a
#       few EMUs when the picture was re-touched)
# This script reproduces both observable effects through the PowerPoint
# object model.

$p = $ppt.ActivePresentation

# Slide 4 in this deck is the one that carries the picture named "그림 145"
# (nvPicPr cNvPr id="1044"). Locate it defensively by Id rather than by
# index/name so the script is resilient to any incidental shape reordering.
$s = $p.Slides.Item(4)

$targetId = 1044
$pic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Id -eq $targetId) {
        $pic = $shp
    }
}

if ($pic -ne $null) {
    # Record where the temp/export copy of the picture lives -- this becomes
    # the shape's "descr" attribute (Alt Text) in the OOXML.
    $pic.AlternativeText = "C:/Users/Admin1/AppData/Roaming/PolarisOffice/ETemp/8536_15965928/fImage688723905436.png"

    # Nudge the picture's position & size (values are in points; PowerPoint
    # stores geometry in EMUs internally -- 1 pt = 12700 EMU).
    #   off  : x 1222375 -> 1230630 EMU , y 3075305 -> 3075305 EMU (unchanged)
    #   ext  : cx 4133850 -> 4126230 EMU, cy 1805940 -> 1779905 EMU
    $pic.Left   = 96.9
    $pic.Top    = 242.15001
    $pic.Width  = 324.90001
    $pic.Height = 140.15001
}
